$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-tracking entry: 2025-09-12, 75 min Projektarbeit, 15 min Unterricht.
# Copy the date-format style from A3 (09/08) onto A4 so the new date cell
# reuses the existing "short date" cell style instead of creating a new one.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("A4").Value = 45912
$ws.Range("B4").Value = 75
$ws.Range("D4").Value = 15

# Matches the saved selection state in the workbook after the edit.
$ws.Range("C10").Select()
